# Applies the commit's highlight-color changes and moves the "_GoBack"
# bookmark from the "Test Data" bullet to the "Website Screenshots" bullet.

$d = $word.ActiveDocument

function Normalize-Highlight($h) {
    # Range.HighlightColorIndex sometimes comes back as the WdColorIndex
    # number and sometimes as the raw OOXML <w:highlight> string (e.g.
    # "magenta" has no WdColorIndex constant). Normalize to the lowercase
    # OOXML highlight name so comparisons are consistent either way.
    if ($h -eq 1) { return "black" }
    if ($h -eq 2) { return "blue" }
    if ($h -eq 3) { return "turquoise" }
    if ($h -eq 4) { return "brightgreen" }
    if ($h -eq 5) { return "pink" }
    if ($h -eq 6) { return "red" }
    if ($h -eq 7) { return "yellow" }
    if ($h -eq 8) { return "white" }
    if ($h -eq 9) { return "darkblue" }
    if ($h -eq 10) { return "teal" }
    if ($h -eq 11) { return "green" }
    if ($h -eq 12) { return "violet" }
    if ($h -eq 13) { return "darkred" }
    if ($h -eq 14) { return "darkyellow" }
    if ($h -eq 15) { return "gray50" }
    if ($h -eq 16) { return "gray25" }
    return ([string]$h).ToLower()
}

function Set-BulletHighlight($paragraphText, $fromColor, $toColor) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd()
        if ($t -eq $paragraphText) {
            $current = Normalize-Highlight $p.Range.HighlightColorIndex
            if ($current -eq $fromColor) {
                $xml = $p.Range.WordOpenXML
                $newXml = $xml.Replace('w:val="' + $fromColor + '"', 'w:val="' + $toColor + '"')
                $p.Range.InsertXML($newXml)
                return $true
            }
        }
    }
    return $false
}

# Logbook: yellow -> green
Set-BulletHighlight "Logbook" "yellow" "green"

# Client Communication: yellow -> green
Set-BulletHighlight "Client Communication" "yellow" "green"

# Prototype Screenshots: magenta -> yellow
Set-BulletHighlight "Prototype Screenshots" "magenta" "yellow"

# Website Screenshots: red -> magenta
Set-BulletHighlight "Website Screenshots" "red" "magenta"

# Move the "_GoBack" bookmark from "Test Data" down to "Website Screenshots".
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "Website Screenshots") {
        $current = Normalize-Highlight $p.Range.HighlightColorIndex
        if ($current -eq "magenta") {
            $target = $d.Range($p.Range.Start, $p.Range.Start)
            $d.Bookmarks.Add("_GoBack", $target)
            break
        }
    }
}
